$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("hospital")

# Update the hospital name (dropdown) and report-by team name labels
$ws.Range("B2").Value = "Siem Reap Provincial Referral Hospital"
$ws.Range("B5").Value = "Battambang Microbiology Team"

# Activate the sheet and move the selection to B2 (matches saved sheetView selection)
$ws.Activate()
$ws.Range("B2").Select()
